{"js": "const newTexts = [\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 10.09.24: \u26a1\ufe0f\ud83d\ude80\",\n  \"Can LLMs Generate Novel Research Ideas? A Large-Scale Human Study with 100+ NLP Researchers\",\n  \"\u05d4\u05d0\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4  \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d9\u05d9\u05e6\u05e8 \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05de\u05d7\u05e7\u05e8 \u05d7\u05d3\u05e9\u05e0\u05d9\u05d9\u05dd? \ud83e\udd14 \u05de\u05d7\u05e7\u05e8 \u05d7\u05d3\u05e9 \u05de\u05e2\u05d5\u05e8\u05e8 \u05d2\u05dc\u05d9\u05dd. \u05e8\u05d0\u05d9\u05e0\u05d5 \u05dc\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05d4\u05ea\u05dc\u05d4\u05d1\u05d5\u05ea \u05e8\u05d1\u05d4 \u05e1\u05d1\u05d9\u05d1 \u05d4\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1-LLMs \u05dc\u05d2\u05d9\u05dc\u05d5\u05d9\u05d9\u05dd \u05de\u05d3\u05e2\u05d9\u05d9\u05dd. \u05d0\u05da \u05d4\u05d0\u05dd \u05d4\u05dd \u05d1\u05d0\u05de\u05ea \u05de\u05e1\u05d5\u05d2\u05dc\u05d9\u05dd \u05dc\u05d4\u05d2\u05d9\u05e2 \u05dc\u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05d7\u05d3\u05e9\u05e0\u05d9\u05d9\u05dd \u05d1\u05e8\u05de\u05ea \u05e8\u05d0\u05d5\u05d9\u05d4 \u05dc\u05d7\u05d5\u05e7\u05e8 \u05d1\u05de\u05d5\u05e1\u05d3 \u05d0\u05e7\u05d3\u05de\u05d9 \u05d0\u05d5 \u05d1\u05ea\u05e2\u05e9\u05d9\u05d4?\",\n  \"\u05de\u05d7\u05d1\u05e8\u05d9 \u05d4\u05de\u05d0\u05de\u05e8 \u05ea\u05db\u05e0\u05e0\u05d5 \u05e0\u05d9\u05e1\u05d5\u05d9 \u05db\u05d3\u05d9 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05ea \u05d4\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4. \u05d4\u05dd \u05e9\u05db\u05e8\u05d5 \u05de\u05e2\u05dc 100 \u05de\u05d5\u05de\u05d7\u05d9 \u05e2\u05d9\u05d1\u05d5\u05d3 \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea \u05dc\u05db\u05ea\u05d5\u05d1 \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05de\u05d7\u05e7\u05e8 \u05d5\u05dc\u05d1\u05d7\u05d5\u05df \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05e9\u05e0\u05d5\u05e6\u05e8\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d1\u05e0\u05d9 \u05d0\u05d3\u05dd \u05d5-LLMs (\u05d1\u05e2\u05d9\u05d5\u05d5\u05e8 \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d1\u05d5\u05d3\u05e7\u05d9\u05dd \u05dc\u05d0 \u05d9\u05d3\u05e2\u05d5 \u05de\u05d4 \u05de\u05e7\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05e9\u05d4\u05dd \u05d1\u05d5\u05d3\u05e7\u05d9\u05dd).\",\n  \"\u05de\u05ea\u05d1\u05e8\u05e8 \u05db\u05d9 \u05d4\u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05d4-LLM \u05e0\u05e9\u05e4\u05d8\u05d5 (\u05d1\u05d0\u05d5\u05e4\u05df \u05dc\u05d0 \u05de\u05e4\u05ea\u05d9\u05e2 \u05e7\u05dc\u05d5\u05d3 \u05e0\u05d1\u05d7\u05e8 \u05dc\u05de\u05e9\u05d9\u05de\u05d4 \u05d6\u05d5) \u05db\u05d7\u05d3\u05e9\u05e0\u05d9\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05de\u05d5\u05de\u05d7\u05d9\u05dd \u05d0\u05e0\u05d5\u05e9\u05d9\u05d9\u05dd (\u05e2\u05dd \u05de\u05d5\u05d1\u05d4\u05e7\u05d5\u05ea \u05e1\u05d8\u05d8\u05d9\u05e1\u05d8\u05d9\u05ea), \u05d0\u05da \u05d3\u05d5\u05e8\u05d2\u05d5 \u05e0\u05de\u05d5\u05da \u05d9\u05d5\u05ea\u05e8 \u05d1\u05d4\u05d9\u05ea\u05db\u05e0\u05d5\u05ea.\",\n  \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05d9\u05e0\u05d9\u05dd \u05d0\u05ea \u05de\u05d4\u05d7\u05d5\u05d6\u05e7\u05d5\u05ea \u05d4\u05d1\u05d0\u05d5\u05ea \u05e9\u05dc \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05d4-LLM:\",\n  \"- \u05d4\u05e6\u05e2\u05ea \u05de\u05db\u05d9\u05dc\u05d4 \u05e9\u05d9\u05dc\u05d5\u05d1\u05d9\u05dd \u05d9\u05d9\u05d7\u05d5\u05d3\u05d9\u05d9\u05dd \u05e9\u05dc \u05d8\u05db\u05e0\u05d9\u05e7\u05d5\u05ea \u05de\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd\",\n  \"- \u05d7\u05e7\u05d9\u05e8\u05ea \u05ea\u05d7\u05d5\u05de\u05d9\u05dd \u05e9\u05dc\u05d0 \u05e0\u05d7\u05e7\u05e8\u05d5 \u05de\u05e1\u05e4\u05d9\u05e7\",\n  \"- \u05d9\u05e6\u05d9\u05e8\u05ea \u05e0\u05d9\u05e1\u05d5\u05d9\u05d9 \u05de\u05d7\u05e9\u05d1\u05d4 \u05d9\u05e6\u05d9\u05e8\u05ea\u05d9\u05d9\u05dd \u05d5\u05de\u05e7\u05d5\u05e8\u05d9\u05d9\u05dd\",\n  \"\u05e2\u05dd \u05d6\u05d0\u05ea, \u05d4\u05d9\u05d5 \u05dc\u05d4\u05dd \u05d2\u05dd \u05db\u05de\u05d4 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05e2\u05d9\u05d5\u05ea\u05d9\u05d5\u05ea:\",\n  \"- \u05d7\u05d5\u05e1\u05e8 \u05e4\u05d9\u05e8\u05d5\u05d8 \u05de\u05e1\u05e4\u05e7 \u05d1\u05e0\u05d5\u05d2\u05e2 \u05dc\u05d9\u05d9\u05e9\u05d5\u05dd\",\n  \"- \u05e9\u05d9\u05de\u05d5\u05e9 \u05dc\u05d0 \u05e0\u05db\u05d5\u05df \u05d1\u05de\u05d0\u05d2\u05e8\u05d9 \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd\",\n  \"- \u05d4\u05d7\u05de\u05e6\u05ea \u05d1\u05d9\u05d9\u05e1\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd (\u05dc\u05d0 \u05de\u05e4\u05ea\u05d9\u05e2 \u05db\u05dc\u05dc)\",\n  \"- \u05d4\u05e0\u05d7\u05d5\u05ea \u05dc\u05d0 \u05de\u05e6\u05d9\u05d0\u05d5\u05ea\u05d9\u05d5\u05ea\",\n  \"\u05dc\u05e2\u05d5\u05de\u05ea \u05d6\u05d0\u05ea, \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d9\u05dd \u05e0\u05d8\u05d5 \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e2\u05d5\u05d2\u05e0\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05d1\u05de\u05d7\u05e7\u05e8 \u05e7\u05d9\u05d9\u05dd \u05d5\u05d1\u05e9\u05d9\u05e7\u05d5\u05dc\u05d9\u05dd \u05de\u05e2\u05e9\u05d9\u05d9\u05dd, \u05d0\u05da \u05dc\u05e2\u05ea\u05d9\u05dd \u05e7\u05e8\u05d5\u05d1\u05d5\u05ea \u05d4\u05d9\u05d5 \u05e4\u05d7\u05d5\u05ea \u05d7\u05d3\u05e9\u05e0\u05d9\u05d9\u05dd, \u05d5\u05d1\u05e0\u05d5 \u05d1\u05d0\u05d5\u05e4\u05df \u05d4\u05d3\u05e8\u05d2\u05ea\u05d9 \u05e2\u05dc \u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d9\u05e6\u05d9\u05d5\u05ea \u05d5\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d9\u05d3\u05d5\u05e2\u05d5\u05ea.\",\n  \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05db\u05d9\u05e8\u05d9\u05dd \u05d1\u05e7\u05d5\u05e9\u05d9 \u05dc\u05e9\u05e4\u05d5\u05d8 \u05d7\u05d3\u05e9\u05e0\u05d5\u05ea, \u05d0\u05e4\u05d9\u05dc\u05d5 \u05e2\u05d1\u05d5\u05e8 \u05de\u05d5\u05de\u05d7\u05d9\u05dd. \u05db\u05e6\u05e2\u05d3 \u05d4\u05d1\u05d0, \u05d4\u05dd \u05d4\u05e6\u05d9\u05e2\u05d5 \u05dc\u05ea\u05ea \u05dc\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05dc\u05de\u05de\u05e9 \u05d0\u05ea \u05d4\u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05d4\u05dc\u05dc\u05d5, \u05db\u05d3\u05d9 \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05dd \u05d3\u05d9\u05e8\u05d5\u05d2\u05d9 \u05d4\u05d7\u05d3\u05e9\u05e0\u05d5\u05ea \u05d5\u05d4\u05d4\u05d9\u05ea\u05db\u05e0\u05d5\u05ea \u05de\u05ea\u05d5\u05e8\u05d2\u05de\u05d9\u05dd \u05dc\u05d4\u05d1\u05d3\u05dc\u05d9\u05dd \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd \u05d1\u05de\u05e6\u05d9\u05d0\u05d5\u05ea.\",\n  \"https://arxiv.org/abs/2409.04109\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items');\nawait context.sync();\n\n// First 9 paragraphs (index 0-8) map 1:1 in place: replace their text.\nconst directCount = 9;\nfor (let i = 0; i < directCount; i++) {\n  paragraphs.items[i].insertText(newTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Paragraph index 9 (old '\u05d1\u05e7\u05d9\u05e6\u05d5\u05e8...') becomes the new '\u05e2\u05dd \u05d6\u05d0\u05ea...' line,\n// then 6 brand-new paragraphs are inserted after it (before the old last/link paragraph),\n// and finally the old last paragraph (the link) gets its text replaced too.\nparagraphs.items[9].insertText(newTexts[9], Word.InsertLocation.replace);\nawait context.sync();\n\nlet anchor = paragraphs.items[9];\nconst insertedTexts = newTexts.slice(10, 16); // 6 new paragraphs\nfor (const t of insertedTexts) {\n  anchor = anchor.insertParagraph(t, Word.InsertLocation.after);\n  await context.sync();\n}\n\n// Last original paragraph (the link) \u2014 now shifted after the inserts \u2014 gets new text.\nparagraphs.items[10].insertText(newTexts[16], Word.InsertLocation.replace);\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$newTexts = @(\n  \"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 10.09.24: \u26a1\ufe0f\ud83d\ude80\",\n  \"Can LLMs Generate Novel Research Ideas? A Large-Scale Human Study with 100+ NLP Researchers\",\n  \"\u05d4\u05d0\u05dd \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4  \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05d9\u05d9\u05e6\u05e8 \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05de\u05d7\u05e7\u05e8 \u05d7\u05d3\u05e9\u05e0\u05d9\u05d9\u05dd? \ud83e\udd14 \u05de\u05d7\u05e7\u05e8 \u05d7\u05d3\u05e9 \u05de\u05e2\u05d5\u05e8\u05e8 \u05d2\u05dc\u05d9\u05dd. \u05e8\u05d0\u05d9\u05e0\u05d5 \u05dc\u05d0\u05d7\u05e8\u05d5\u05e0\u05d4 \u05d4\u05ea\u05dc\u05d4\u05d1\u05d5\u05ea \u05e8\u05d1\u05d4 \u05e1\u05d1\u05d9\u05d1 \u05d4\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1-LLMs \u05dc\u05d2\u05d9\u05dc\u05d5\u05d9\u05d9\u05dd \u05de\u05d3\u05e2\u05d9\u05d9\u05dd. \u05d0\u05da \u05d4\u05d0\u05dd \u05d4\u05dd \u05d1\u05d0\u05de\u05ea \u05de\u05e1\u05d5\u05d2\u05dc\u05d9\u05dd \u05dc\u05d4\u05d2\u05d9\u05e2 \u05dc\u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05d7\u05d3\u05e9\u05e0\u05d9\u05d9\u05dd \u05d1\u05e8\u05de\u05ea \u05e8\u05d0\u05d5\u05d9\u05d4 \u05dc\u05d7\u05d5\u05e7\u05e8 \u05d1\u05de\u05d5\u05e1\u05d3 \u05d0\u05e7\u05d3\u05de\u05d9 \u05d0\u05d5 \u05d1\u05ea\u05e2\u05e9\u05d9\u05d4?\",\n  \"\u05de\u05d7\u05d1\u05e8\u05d9 \u05d4\u05de\u05d0\u05de\u05e8 \u05ea\u05db\u05e0\u05e0\u05d5 \u05e0\u05d9\u05e1\u05d5\u05d9 \u05db\u05d3\u05d9 \u05dc\u05d1\u05d3\u05d5\u05e7 \u05d0\u05ea \u05d4\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4. \u05d4\u05dd \u05e9\u05db\u05e8\u05d5 \u05de\u05e2\u05dc 100 \u05de\u05d5\u05de\u05d7\u05d9 \u05e2\u05d9\u05d1\u05d5\u05d3 \u05e9\u05e4\u05d4 \u05d8\u05d1\u05e2\u05d9\u05ea \u05dc\u05db\u05ea\u05d5\u05d1 \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05de\u05d7\u05e7\u05e8 \u05d5\u05dc\u05d1\u05d7\u05d5\u05df \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05e9\u05e0\u05d5\u05e6\u05e8\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d1\u05e0\u05d9 \u05d0\u05d3\u05dd \u05d5-LLMs (\u05d1\u05e2\u05d9\u05d5\u05d5\u05e8 \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d1\u05d5\u05d3\u05e7\u05d9\u05dd \u05dc\u05d0 \u05d9\u05d3\u05e2\u05d5 \u05de\u05d4 \u05de\u05e7\u05d5\u05e8\u05d4 \u05e9\u05dc \u05d4\u05e8\u05e2\u05d9\u05d5\u05df \u05e9\u05d4\u05dd \u05d1\u05d5\u05d3\u05e7\u05d9\u05dd).\",\n  \"\u05de\u05ea\u05d1\u05e8\u05e8 \u05db\u05d9 \u05d4\u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05d4-LLM \u05e0\u05e9\u05e4\u05d8\u05d5 (\u05d1\u05d0\u05d5\u05e4\u05df \u05dc\u05d0 \u05de\u05e4\u05ea\u05d9\u05e2 \u05e7\u05dc\u05d5\u05d3 \u05e0\u05d1\u05d7\u05e8 \u05dc\u05de\u05e9\u05d9\u05de\u05d4 \u05d6\u05d5) \u05db\u05d7\u05d3\u05e9\u05e0\u05d9\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05de\u05d5\u05de\u05d7\u05d9\u05dd \u05d0\u05e0\u05d5\u05e9\u05d9\u05d9\u05dd (\u05e2\u05dd \u05de\u05d5\u05d1\u05d4\u05e7\u05d5\u05ea \u05e1\u05d8\u05d8\u05d9\u05e1\u05d8\u05d9\u05ea), \u05d0\u05da \u05d3\u05d5\u05e8\u05d2\u05d5 \u05e0\u05de\u05d5\u05da \u05d9\u05d5\u05ea\u05e8 \u05d1\u05d4\u05d9\u05ea\u05db\u05e0\u05d5\u05ea.\",\n  \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05d9\u05e0\u05d9\u05dd \u05d0\u05ea \u05de\u05d4\u05d7\u05d5\u05d6\u05e7\u05d5\u05ea \u05d4\u05d1\u05d0\u05d5\u05ea \u05e9\u05dc \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05d4-LLM:\",\n  \"- \u05d4\u05e6\u05e2\u05ea \u05de\u05db\u05d9\u05dc\u05d4 \u05e9\u05d9\u05dc\u05d5\u05d1\u05d9\u05dd \u05d9\u05d9\u05d7\u05d5\u05d3\u05d9\u05d9\u05dd \u05e9\u05dc \u05d8\u05db\u05e0\u05d9\u05e7\u05d5\u05ea \u05de\u05d3\u05d5\u05de\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d5\u05e0\u05d9\u05dd\",\n  \"- \u05d7\u05e7\u05d9\u05e8\u05ea \u05ea\u05d7\u05d5\u05de\u05d9\u05dd \u05e9\u05dc\u05d0 \u05e0\u05d7\u05e7\u05e8\u05d5 \u05de\u05e1\u05e4\u05d9\u05e7\",\n  \"- \u05d9\u05e6\u05d9\u05e8\u05ea \u05e0\u05d9\u05e1\u05d5\u05d9\u05d9 \u05de\u05d7\u05e9\u05d1\u05d4 \u05d9\u05e6\u05d9\u05e8\u05ea\u05d9\u05d9\u05dd \u05d5\u05de\u05e7\u05d5\u05e8\u05d9\u05d9\u05dd\",\n  \"\u05e2\u05dd \u05d6\u05d0\u05ea, \u05d4\u05d9\u05d5 \u05dc\u05d4\u05dd \u05d2\u05dd \u05db\u05de\u05d4 \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05d1\u05e2\u05d9\u05d5\u05ea\u05d9\u05d5\u05ea:\",\n  \"- \u05d7\u05d5\u05e1\u05e8 \u05e4\u05d9\u05e8\u05d5\u05d8 \u05de\u05e1\u05e4\u05e7 \u05d1\u05e0\u05d5\u05d2\u05e2 \u05dc\u05d9\u05d9\u05e9\u05d5\u05dd\",\n  \"- \u05e9\u05d9\u05de\u05d5\u05e9 \u05dc\u05d0 \u05e0\u05db\u05d5\u05df \u05d1\u05de\u05d0\u05d2\u05e8\u05d9 \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd\",\n  \"- \u05d4\u05d7\u05de\u05e6\u05ea \u05d1\u05d9\u05d9\u05e1\u05dc\u05d9\u05d9\u05e0\u05d9\u05dd (\u05dc\u05d0 \u05de\u05e4\u05ea\u05d9\u05e2 \u05db\u05dc\u05dc)\",\n  \"- \u05d4\u05e0\u05d7\u05d5\u05ea \u05dc\u05d0 \u05de\u05e6\u05d9\u05d0\u05d5\u05ea\u05d9\u05d5\u05ea\",\n  \"\u05dc\u05e2\u05d5\u05de\u05ea \u05d6\u05d0\u05ea, \u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d9\u05dd \u05e0\u05d8\u05d5 \u05dc\u05d4\u05d9\u05d5\u05ea \u05de\u05e2\u05d5\u05d2\u05e0\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05d1\u05de\u05d7\u05e7\u05e8 \u05e7\u05d9\u05d9\u05dd \u05d5\u05d1\u05e9\u05d9\u05e7\u05d5\u05dc\u05d9\u05dd \u05de\u05e2\u05e9\u05d9\u05d9\u05dd, \u05d0\u05da \u05dc\u05e2\u05ea\u05d9\u05dd \u05e7\u05e8\u05d5\u05d1\u05d5\u05ea \u05d4\u05d9\u05d5 \u05e4\u05d7\u05d5\u05ea \u05d7\u05d3\u05e9\u05e0\u05d9\u05d9\u05dd, \u05d5\u05d1\u05e0\u05d5 \u05d1\u05d0\u05d5\u05e4\u05df \u05d4\u05d3\u05e8\u05d2\u05ea\u05d9 \u05e2\u05dc \u05d0\u05d9\u05e0\u05d8\u05d5\u05d0\u05d9\u05e6\u05d9\u05d5\u05ea \u05d5\u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d9\u05d3\u05d5\u05e2\u05d5\u05ea.\",\n  \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05d9\u05e0\u05d9\u05dd \u05e9\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05db\u05d9\u05e8\u05d9\u05dd \u05d1\u05e7\u05d5\u05e9\u05d9 \u05dc\u05e9\u05e4\u05d5\u05d8 \u05d7\u05d3\u05e9\u05e0\u05d5\u05ea, \u05d0\u05e4\u05d9\u05dc\u05d5 \u05e2\u05d1\u05d5\u05e8 \u05de\u05d5\u05de\u05d7\u05d9\u05dd. \u05db\u05e6\u05e2\u05d3 \u05d4\u05d1\u05d0, \u05d4\u05dd \u05d4\u05e6\u05d9\u05e2\u05d5 \u05dc\u05ea\u05ea \u05dc\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05dc\u05de\u05de\u05e9 \u05d0\u05ea \u05d4\u05e8\u05e2\u05d9\u05d5\u05e0\u05d5\u05ea \u05d4\u05dc\u05dc\u05d5, \u05db\u05d3\u05d9 \u05dc\u05e8\u05d0\u05d5\u05ea \u05d0\u05dd \u05d3\u05d9\u05e8\u05d5\u05d2\u05d9 \u05d4\u05d7\u05d3\u05e9\u05e0\u05d5\u05ea \u05d5\u05d4\u05d4\u05d9\u05ea\u05db\u05e0\u05d5\u05ea \u05de\u05ea\u05d5\u05e8\u05d2\u05de\u05d9\u05dd \u05dc\u05d4\u05d1\u05d3\u05dc\u05d9\u05dd \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05d9\u05dd \u05d1\u05de\u05e6\u05d9\u05d0\u05d5\u05ea.\",\n  \"https://arxiv.org/abs/2409.04109\",\n)\n\n# Paragraphs 1-9 (1-indexed) map directly onto the new text in place.\nfor ($i = 1; $i -le 9; $i++) {\n  $d.Paragraphs($i).Range.Text = $newTexts[$i - 1]\n}\n\n# Paragraph 10 (old closing line) becomes the new '\u05e2\u05dd \u05d6\u05d0\u05ea...' transition line.\n$d.Paragraphs(10).Range.Text = $newTexts[9]\n\n# Insert 6 brand-new paragraphs after paragraph 10, before the old last (link) paragraph.\n$insertPos = 10\nfor ($i = 10; $i -le 15; $i++) {\n  $d.Paragraphs($insertPos).Range.InsertParagraphAfter()\n  $insertPos = $insertPos + 1\n  $d.Paragraphs($insertPos).Range.Text = $newTexts[$i]\n}\n\n# Final paragraph is the old last (link) paragraph, now pushed down to position 17; replace its text.\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = $newTexts[16]"}
